# "Update countries & provincias Spain"
#
# The "Pais" sheet is a COVID-19 table sorted descending by column B
# ("Casos totales"). This refreshes the statistics (columns B:H) for the
# countries whose figures changed in the new data pull, and bumps the
# "last updated" timestamp. Column A (country name reference) is never
# touched on any row -- a few countries' updated totals now rank higher
# than their neighbours, so those neighbouring rows' B:H values shift down
# one row to keep the table sorted, without moving the country labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp: 18:20 -> 18:50
$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 18:50"

# Estados Unidos
$ws.Range("B4").Value = 200289
$ws.Range("C4").Value = 11759
$ws.Range("D4").Value = 8707
$ws.Range("E4").Value = 187188
$ws.Range("F4").Value = 4888
$ws.Range("G4").Value = 341
$ws.Range("H4").Value = 4394

# Irlanda
$ws.Range("B26").Value = 3447
$ws.Range("C26").Value = 212
$ws.Range("E26").Value = 3357
$ws.Range("G26").Value = 14
$ws.Range("H26").Value = 85

# Ecuador's total now overtakes Rumania/Polonia/Luxemburgo/Filipinas, so it
# moves to the top of this block (row 31); those four countries' data
# shifts down one row each (32-35) keeping the column-B sort intact.
$ws.Range("B31").Value = 2748
$ws.Range("C31").Value = 446
$ws.Range("D31").Value = 58
$ws.Range("E31").Value = 2597
$ws.Range("F31").Value = 100
$ws.Range("G31").Value = 14
$ws.Range("H31").Value = 93

$ws.Range("B32").Value = 2460
$ws.Range("C32").Value = 215
$ws.Range("D32").Value = 252
$ws.Range("E32").Value = 2122
$ws.Range("F32").Value = 57
$ws.Range("G32").Value = 4
$ws.Range("H32").Value = 86

$ws.Range("B33").Value = 2420
$ws.Range("C33").Value = 109
$ws.Range("D33").Value = 7
$ws.Range("E33").Value = 2377
$ws.Range("F33").Value = 50
$ws.Range("G33").Value = 3
$ws.Range("H33").Value = 36

$ws.Range("B34").Value = 2319
$ws.Range("C34").Value = 141
$ws.Range("D34").Value = 80
$ws.Range("E34").Value = 2210
$ws.Range("F34").Value = 31
$ws.Range("G34").Value = 6
$ws.Range("H34").Value = 29

$ws.Range("B35").Value = 2311
$ws.Range("C35").Value = 227
$ws.Range("D35").Value = 50
$ws.Range("F35").Value = 1
$ws.Range("G35").Value = 8
$ws.Range("H35").Value = 96

# Sudafrica
$ws.Range("B44").Value = 1380
$ws.Range("C44").Value = 27
$ws.Range("E44").Value = 1325

# Argelia
$ws.Range("B55").Value = 847
$ws.Range("C55").Value = 131
$ws.Range("D55").Value = 61
$ws.Range("E55").Value = 728

# Armenia
$ws.Range("E68").Value = 536
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = 4

# Barein
$ws.Range("D69").Value = 337
$ws.Range("E69").Value = 226

# Moldavia
$ws.Range("D74").Value = 23
$ws.Range("E74").Value = 395

# Mauricio
$ws.Range("B105").Value = 158
$ws.Range("C105").Value = 15
$ws.Range("E105").Value = 153

# Montenegro
$ws.Range("B112").Value = 123
$ws.Range("C112").Value = 14
$ws.Range("E112").Value = 121

# Monaco's total now overtakes Banglades, so it moves up to row 129;
# Banglades' data (unchanged from the previous pull) shifts down to row 130.
$ws.Range("B129").Value = 55
$ws.Range("D129").Value = 2
$ws.Range("E129").Value = 52
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 1

$ws.Range("B130").Value = 54
$ws.Range("C130").Value = 3
$ws.Range("D130").Value = 25
$ws.Range("E130").Value = 23
$ws.Range("F130").Value = 1
$ws.Range("G130").Value = 1
$ws.Range("H130").Value = 6

# Sudan
$ws.Range("D185").Value = 2
$ws.Range("E185").Value = 3
